$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.220779895782471
$ws.Range("B1").Value = 2.47726583480835
$ws.Range("C1").Value = 7.358223915100098
$ws.Range("D1").Value = 2.23894476890564
$ws.Range("E1").Value = 1.159412145614624
